# ---------------------------------------------------------------------------
# Team07Report.xlsx edit script
# "updated burndown sheet, backlog and sprint one"
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Backlog sheet: mark the three "Coding" rows as "Done"
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("E2").Value = "Done"
$backlog.Range("E4").Value = "Done"
$backlog.Range("E9").Value = "Done"

# ---------------------------------------------------------------------------
# 2) Sprint1 sheet: mark the three "Coding" rows as "Done", add actual
#    size/time + "Yes" completed flag, and add the two owner cells that
#    were missing on the "done" subtasks.
# ---------------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

$sprint1.Range("D2").Value = "Done"
$sprint1.Range("G2").Value = 5
$sprint1.Range("H2").Value = 5
$sprint1.Range("I2").Value = "Yes"

$sprint1.Range("C4").Value = "rp"

$sprint1.Range("D6").Value = "Done"
$sprint1.Range("G6").Value = 5
$sprint1.Range("H6").Value = 15
$sprint1.Range("I6").Value = "Yes"

$sprint1.Range("D10").Value = "Done"
$sprint1.Range("G10").Value = 50
$sprint1.Range("H10").Value = 20
$sprint1.Range("I10").Value = "Yes"

$sprint1.Range("C11").Value = "rp"

# ---------------------------------------------------------------------------
# 3) Burndown sheet: insert a "Sprint" label column in front of the table
#    and add the second data row (sprint 1 results).
# ---------------------------------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Columns("A:A").Insert()

$burndown.Range("A1").Value = "Sprint"
$burndown.Range("A2").Value = "Start"
$burndown.Range("A3").Value = "Sprint 1"
$burndown.Range("A4").Value = "Sprint 2"
$burndown.Range("A5").Value = "Sprint 3"
$burndown.Range("A6").Value = "Sprint 4"

$burndown.Range("B3").Value = 42810
$burndown.Range("C3").Value = 33
$burndown.Range("E3").Value = 60
$burndown.Range("F3").Value = 40
$burndown.Range("G3").Formula = "=E3/F3"

# Move the burndown chart one column to the right so it stays clear of the
# new "Sprint" label column.
$burndownChart = $burndown.ChartObjects().Item(1)
$colAWidth = $burndown.Columns.Item(1).Width()
$burndownChart.Left = $burndownChart.Left() + $colAWidth

# ---------------------------------------------------------------------------
# 4) View-state touch-ups
# ---------------------------------------------------------------------------

# Stories: scroll back to the top-left (used to be parked at A33).
$stories = $wb.Worksheets.Item("Stories")
$stories.Activate()
$win = $wb.Windows().Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$stories.Range("A39:B39").Select()

# Sprint1: scroll back to the top-left and move the selection.
$sprint1.Activate()
$win.ScrollRow = 1
$win.ScrollColumn = 1
$sprint1.Range("H6").Select()

# Burndown: move the selection to reflect the new layout.
$burndown.Activate()
$burndown.Range("G4").Select()

# Team (E9 region) / Backlog selection unaffected; finish on Backlog, which
# becomes the active tab.
$backlog.Activate()
$backlog.Range("E9").Select()
